$wb = $excel.ActiveWorkbook
$table = $wb.Worksheets.Item(1)

# --- Bump the internal sheetId counter so the new sheet lands on sheetId=3 ---
# (mirrors the original commit, where the "Submissions" sheet ended up with
# sheetId="3" even though it is only the 2nd sheet / r:id="rId2")
$bump = $wb.Worksheets.Add($null, $table)
$bump.Name = "Temp1"

# --- Duplicate "Table" right after itself; the duplicate becomes "Submissions" ---
$table.Copy($null, $bump)
$subs = $wb.Worksheets.Item(3)
$subs.Name = "Submissions"

# Remove the scratch sheet used only to advance the sheetId counter
$bump.Delete()

# Re-select/activate the Submissions sheet (now at position 2) so it becomes
# the active tab, matching activeTab="1" / tabSelected="1" in the target file
$wb.Worksheets.Item(2).Activate()
$subs = $wb.Worksheets.Item(2)

# --- Fill in the Submissions table content ---
# Header row
$subs.Range("A1").Value = "Attribute"
$subs.Range("B1").Value = "Datatype"

# USERID / VARCHAR(15)
$subs.Range("A2").Value = "USERID"
$subs.Range("B2").Value = "VARCHAR(15)"

# SUBMISSIONID(PK) / NUMBER  (datatype entered before the attribute name,
# same insertion order as the original workbook's shared-string table)
$subs.Range("B3").Value = "NUMBER"
$subs.Range("A3").Value = "SUBMISSIONID(PK)"

# RATING / NUMBER
$subs.Range("A4").Value = "RATING"
$subs.Range("B4").Value = "NUMBER"

# SUBMISSION_DATE / DATE (again datatype entered first)
$subs.Range("B5").Value = "DATE"
$subs.Range("A5").Value = "SUBMISSION_DATE"

# --- Column widths for the new sheet ---
$subs.Columns.Item(1).ColumnWidth = 18.1
$subs.Columns.Item(2).ColumnWidth = 12.33
